$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (Nigeria / "2")
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 0.17105
$ws.Range("E2").Value = 0.244
$ws.Range("G2").Value = -0.003733876442634084
$ws.Range("H2").Value = -0.003733876442634084
$ws.Range("I2").Value = -0.02783435166327223
$ws.Range("J2").Value = -0.02783435166327223
$ws.Range("K2").Value = -5.099999999999998
$ws.Range("L2").Value = -0.01731160896130345
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = -0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 70.48999999999999
$ws.Range("V2").Value = 0.9790277777777777
$ws.Range("W2").Value = 0.8692423105776443
$ws.Range("X2").Value = 0.1121406688681119
$ws.Range("Y2").Value = 0.7571016417095324
$ws.Range("Z2").Value = 19.66622162883845
$ws.Range("AA2").Value = 0.8278485467910276
$ws.Range("AB2").Value = 0.1074833294770007
$ws.Range("AC2").Value = 0.7203652173140268
$ws.Range("AD2").Value = 5.12
$ws.Range("AF2").Value = 5.12
$ws.Range("AG2").Value = -65.36999999999999
$ws.Range("AH2").Value = 0.06639004149377593
$ws.Range("AI2").Value = 0.09188801148600145
$ws.Range("AJ2").Value = -9.859728506787315
$ws.Range("AK2").Value = 4.425863236289777
$ws.Range("AL2").Value = 1.402
$ws.Range("AM2").Value = 1.402
$ws.Range("AN2").Value = -0.8258064516129033
$ws.Range("AO2").Value = -5.848787446504992
$ws.Range("AP2").Value = 10.54354838709677
$ws.Range("AQ2").Value = -5.848787446504992

# ---------------------------------------------------------------------------
# Row 3 (AIICO -> African Alliance Insurance Plc)
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "African Alliance Insurance Plc (NGSE:AFRINSURE)"
$ws.Range("D3").Value = 0.0131
$ws.Range("E3").ClearContents()
$ws.Range("G3").Value = -0.5412621359223301
$ws.Range("H3").Value = -0.5412621359223301
$ws.Range("I3").Value = -0.6140776699029126
$ws.Range("J3").Value = -0.6140776699029126
$ws.Range("K3").Value = -21.9
$ws.Range("L3").Value = -0.5315533980582523
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 2.49
$ws.Range("V3").Value = 0.2305555555555556
$ws.Range("W3").Value = 1.412903225806452
$ws.Range("X3").Value = 0.1136635459207394
$ws.Range("Y3").Value = 1.299239679885712
$ws.Range("Z3").Value = -1.929742388758782
$ws.Range("AA3").Value = 1.185011709601874
$ws.Range("AB3").Value = 0.1078881594477103
$ws.Range("AC3").Value = 1.077123550154163
$ws.Range("AD3").Value = 1.18
$ws.Range("AF3").Value = 1.18
$ws.Range("AG3").Value = -1.31
$ws.Range("AH3").Value = 0.09849749582637729
$ws.Range("AI3").Value = -0.03369503141062249
$ws.Range("AJ3").Value = -0.1380400421496312
$ws.Range("AK3").Value = 0.03492402026126366
$ws.Range("AL3").Value = 0.525
$ws.Range("AM3").Value = 0.525
$ws.Range("AN3").Value = -0.04738955823293173
$ws.Range("AO3").Value = -48.19047619047619
$ws.Range("AP3").Value = 0.05261044176706828
$ws.Range("AQ3").Value = -48.19047619047619

# ---------------------------------------------------------------------------
# Row 4 (African Alliance Insurance Plc -> AIICO)
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "AIICO Insurance Plc (NGSE:AIICO)"
$ws.Range("D4").Value = 0.329
$ws.Range("E4").Value = 0.244
$ws.Range("G4").Value = 0.08366219415943173
$ws.Range("H4").Value = 0.08366219415943173
$ws.Range("I4").Value = 0.06748224151539069
$ws.Range("J4").Value = 0.06748224151539069
$ws.Range("K4").Value = 16.8
$ws.Range("L4").Value = 0.06629834254143646
$ws.Range("U4").Value = 68
$ws.Range("V4").Value = 1.111111111111111
$ws.Range("W4").Value = 0.3255813953488372
$ws.Range("X4").Value = 0.1106177918154844
$ws.Range("Y4").Value = 0.2149636035333529
$ws.Range("Z4").Value = 6.974951830443159
$ws.Range("AA4").Value = 0.4706853839801816
$ws.Range("AB4").Value = 0.1070784995062912
$ws.Range("AC4").Value = 0.3636068844738904
$ws.Range("AD4").Value = 3.94
$ws.Range("AF4").Value = 3.94
$ws.Range("AG4").Value = -64.06
$ws.Range("AH4").Value = 0.0604851089960086
$ws.Range("AI4").Value = 0.04342076261847036
$ws.Range("AJ4").Value = 22.3986013986014
$ws.Range("AK4").Value = -2.817062445030783
$ws.Range("AL4").Value = 0.877
$ws.Range("AM4").Value = 0.877
$ws.Range("AN4").Value = 0.2106951871657754
$ws.Range("AO4").Value = 19.49828962371722
$ws.Range("AP4").Value = -3.425668449197861
$ws.Range("AQ4").Value = 19.49828962371722
